$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13: EPAM Systems, Inc. / 1 / US
# (company name cell uses an 8pt Segoe UI font, color #212529)
$ws.Range("A13").Value = "EPAM Systems, Inc."
$ws.Range("A13").Font.Name = "Segoe UI"
$ws.Range("A13").Font.Size = 8
$ws.Range("A13").Font.Color = 2696481
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "US"

# New row 14: Tech Mahindra Limited / 4 / India
$ws.Range("A14").Value = "Tech Mahindra Limited"
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = "India"

# Match the author's final selection/scroll position
$null = $ws.Range("L9").Select()
